$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 22-31: timestamp (A), label (B)
$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"

# Update sensor columns C:H for rows 2-31 with refreshed data
$ws.Cells.Item(2, 3).Value = -9.514568328857422
$ws.Cells.Item(2, 4).Value = -7.287443161010742
$ws.Cells.Item(2, 5).Value = -4.99766206741333
$ws.Cells.Item(2, 6).Value = 1.033005767512607
$ws.Cells.Item(2, 7).Value = -2.26011000282463
$ws.Cells.Item(2, 8).Value = -0.9051166365201431
$ws.Cells.Item(3, 3).Value = -7.823380947113037
$ws.Cells.Item(3, 4).Value = -6.074180126190186
$ws.Cells.Item(3, 5).Value = -4.752357959747314
$ws.Cells.Item(3, 6).Value = 0.9114278760449622
$ws.Cells.Item(3, 7).Value = -1.322062040197435
$ws.Cells.Item(3, 8).Value = -0.8478403153090641
$ws.Cells.Item(4, 3).Value = -3.637385606765747
$ws.Cells.Item(4, 4).Value = -2.725360631942749
$ws.Cells.Item(4, 5).Value = -3.621543645858765
$ws.Cells.Item(4, 6).Value = 0.6373885839149871
$ws.Cells.Item(4, 7).Value = -0.4370783201750674
$ws.Cells.Item(4, 8).Value = -0.9901816481831425
$ws.Cells.Item(5, 3).Value = -10.94634628295898
$ws.Cells.Item(5, 4).Value = -2.156153917312622
$ws.Cells.Item(5, 5).Value = -9.503169059753418
$ws.Cells.Item(5, 6).Value = -0.08095168998871549
$ws.Cells.Item(5, 7).Value = -0.8572516530409613
$ws.Cells.Item(5, 8).Value = -0.2425374309907035
$ws.Cells.Item(6, 3).Value = 66.69537353515625
$ws.Cells.Item(6, 4).Value = -35.59264755249023
$ws.Cells.Item(6, 5).Value = -8.555927276611328
$ws.Cells.Item(6, 6).Value = -0.5983732620189923
$ws.Cells.Item(6, 7).Value = -0.3525244653353435
$ws.Cells.Item(6, 8).Value = -0.1835469028045379
$ws.Cells.Item(7, 3).Value = -6.669784069061279
$ws.Cells.Item(7, 4).Value = -7.415677547454834
$ws.Cells.Item(7, 5).Value = 10.31031608581543
$ws.Cells.Item(7, 6).Value = 0.904387354850761
$ws.Cells.Item(7, 7).Value = 0.4711575967570125
$ws.Cells.Item(7, 8).Value = -3.501346578200644
$ws.Cells.Item(8, 3).Value = -6.682662963867188
$ws.Cells.Item(8, 4).Value = 0.0761735439300537
$ws.Cells.Item(8, 5).Value = -2.404594898223877
$ws.Cells.Item(8, 6).Value = -0.0358646748394702
$ws.Cells.Item(8, 7).Value = -1.06903420953914
$ws.Cells.Item(8, 8).Value = -1.137895204212497
$ws.Cells.Item(9, 3).Value = -4.13407564163208
$ws.Cells.Item(9, 4).Value = -28.60597419738769
$ws.Cells.Item(9, 5).Value = 8.249073028564453
$ws.Cells.Item(9, 6).Value = -2.645820761549075
$ws.Cells.Item(9, 7).Value = -2.995684344193039
$ws.Cells.Item(9, 8).Value = 2.325931413420326
$ws.Cells.Item(10, 3).Value = -12.32790374755859
$ws.Cells.Item(10, 4).Value = 6.361005783081055
$ws.Cells.Item(10, 5).Value = -18.50937080383301
$ws.Cells.Item(10, 6).Value = -5.634616028303401
$ws.Cells.Item(10, 7).Value = 1.589792383128232
$ws.Cells.Item(10, 8).Value = 0.8637641414828536
$ws.Cells.Item(11, 3).Value = -0.0332281589508056
$ws.Cells.Item(11, 4).Value = -8.605781555175781
$ws.Cells.Item(11, 5).Value = -5.739476203918457
$ws.Cells.Item(11, 6).Value = 0.2434335577077258
$ws.Cells.Item(11, 7).Value = -0.9524030027718666
$ws.Cells.Item(11, 8).Value = -0.7793511700356166
$ws.Cells.Item(12, 3).Value = -14.07493591308594
$ws.Cells.Item(12, 4).Value = -31.83533477783203
$ws.Cells.Item(12, 5).Value = -4.278344631195068
$ws.Cells.Item(12, 6).Value = 3.219139538962262
$ws.Cells.Item(12, 7).Value = -4.921159070113614
$ws.Cells.Item(12, 8).Value = -0.2975026624726805
$ws.Cells.Item(13, 3).Value = 2.169375419616699
$ws.Cells.Item(13, 4).Value = 6.375825881958008
$ws.Cells.Item(13, 5).Value = 16.93547058105469
$ws.Cells.Item(13, 6).Value = 6.111283532504389
$ws.Cells.Item(13, 7).Value = 1.10011664204216
$ws.Cells.Item(13, 8).Value = 1.229119239182312
$ws.Cells.Item(14, 3).Value = -3.887731313705444
$ws.Cells.Item(14, 4).Value = 1.673339605331421
$ws.Cells.Item(14, 5).Value = 2.357208251953125
$ws.Cells.Item(14, 6).Value = 2.061637947614152
$ws.Cells.Item(14, 7).Value = 3.835073033968605
$ws.Cells.Item(14, 8).Value = -0.2914631209154284
$ws.Cells.Item(15, 3).Value = 3.855255126953125
$ws.Cells.Item(15, 4).Value = -33.85980224609375
$ws.Cells.Item(15, 5).Value = 3.360419273376465
$ws.Cells.Item(15, 6).Value = -3.819613575935366
$ws.Cells.Item(15, 7).Value = 2.543305224385755
$ws.Cells.Item(15, 8).Value = 1.029467895113191
$ws.Cells.Item(16, 3).Value = 30.62849044799805
$ws.Cells.Item(16, 4).Value = 7.678761005401611
$ws.Cells.Item(16, 5).Value = -9.237998962402344
$ws.Cells.Item(16, 6).Value = -4.798000733057658
$ws.Cells.Item(16, 7).Value = 7.385336404559252
$ws.Cells.Item(16, 8).Value = 1.159243436466003
$ws.Cells.Item(17, 3).Value = -11.34725379943848
$ws.Cells.Item(17, 4).Value = -16.60527801513672
$ws.Cells.Item(17, 5).Value = -12.14533615112305
$ws.Cells.Item(17, 6).Value = -0.8656106913226824
$ws.Cells.Item(17, 7).Value = 6.255807002385461
$ws.Cells.Item(17, 8).Value = -3.608730400088184
$ws.Cells.Item(18, 3).Value = 6.219323635101318
$ws.Cells.Item(18, 4).Value = -10.72451782226562
$ws.Cells.Item(18, 5).Value = 26.53547286987305
$ws.Cells.Item(18, 6).Value = 4.937833545536877
$ws.Cells.Item(18, 7).Value = -2.039471638613807
$ws.Cells.Item(18, 8).Value = -6.196790481435843
$ws.Cells.Item(19, 3).Value = -4.500537395477295
$ws.Cells.Item(19, 4).Value = 9.60122776031494
$ws.Cells.Item(19, 5).Value = -3.719542026519776
$ws.Cells.Item(19, 6).Value = 5.261263814465728
$ws.Cells.Item(19, 7).Value = -3.646456844505198
$ws.Cells.Item(19, 8).Value = -4.444081427037008
$ws.Cells.Item(20, 3).Value = -59.25642776489258
$ws.Cells.Item(20, 4).Value = -72.75296783447266
$ws.Cells.Item(20, 5).Value = 58.0263671875
$ws.Cells.Item(20, 6).Value = -0.0522782514835689
$ws.Cells.Item(20, 7).Value = 0.8027088176245329
$ws.Cells.Item(20, 8).Value = -0.5144340389076021
$ws.Cells.Item(21, 3).Value = 43.50658416748047
$ws.Cells.Item(21, 4).Value = 8.478635787963867
$ws.Cells.Item(21, 5).Value = -37.43244552612305
$ws.Cells.Item(21, 6).Value = -7.4421002289345
$ws.Cells.Item(21, 7).Value = -0.806786348079851
$ws.Cells.Item(21, 8).Value = 5.891316611191369
$ws.Cells.Item(22, 3).Value = -19.30278778076172
$ws.Cells.Item(22, 4).Value = -6.771676063537598
$ws.Cells.Item(22, 5).Value = -17.75639343261719
$ws.Cells.Item(22, 6).Value = -0.6881247882184418
$ws.Cells.Item(22, 7).Value = -10.85945387133243
$ws.Cells.Item(22, 8).Value = 4.675690663957009
$ws.Cells.Item(23, 3).Value = -18.29881477355957
$ws.Cells.Item(23, 4).Value = -37.5744743347168
$ws.Cells.Item(23, 5).Value = 5.842066764831543
$ws.Cells.Item(23, 6).Value = 4.403560649389491
$ws.Cells.Item(23, 7).Value = -11.45763061786516
$ws.Cells.Item(23, 8).Value = 1.461910155997879
$ws.Cells.Item(24, 3).Value = -2.273155212402344
$ws.Cells.Item(24, 4).Value = 8.59691047668457
$ws.Cells.Item(24, 5).Value = -6.313179969787598
$ws.Cells.Item(24, 6).Value = 4.784720346845424
$ws.Cells.Item(24, 7).Value = 3.472261708358207
$ws.Cells.Item(24, 8).Value = -1.690446103441295
$ws.Cells.Item(25, 3).Value = -0.403256893157959
$ws.Cells.Item(25, 4).Value = 4.687671661376953
$ws.Cells.Item(25, 5).Value = -1.856612205505371
$ws.Cells.Item(25, 6).Value = 0.8494467159797106
$ws.Cells.Item(25, 7).Value = 5.416543818924583
$ws.Cells.Item(25, 8).Value = 0.3924825684777868
$ws.Cells.Item(26, 3).Value = 4.316394805908203
$ws.Cells.Item(26, 4).Value = -26.35572052001953
$ws.Cells.Item(26, 5).Value = -17.98580360412598
$ws.Cells.Item(26, 6).Value = -4.816584700825577
$ws.Cells.Item(26, 7).Value = 1.840564275695681
$ws.Cells.Item(26, 8).Value = 3.20516648785821
$ws.Cells.Item(27, 3).Value = 31.11298370361328
$ws.Cells.Item(27, 4).Value = 3.278896331787109
$ws.Cells.Item(27, 5).Value = -11.80455207824707
$ws.Cells.Item(27, 6).Value = -4.067581341184413
$ws.Cells.Item(27, 7).Value = 7.491072893142713
$ws.Cells.Item(27, 8).Value = -0.0797111165934683
$ws.Cells.Item(28, 3).Value = -10.42159271240234
$ws.Cells.Item(28, 4).Value = -19.3218994140625
$ws.Cells.Item(28, 5).Value = -14.07432746887207
$ws.Cells.Item(28, 6).Value = -0.603338341945886
$ws.Cells.Item(28, 7).Value = 10.07777972605036
$ws.Cells.Item(28, 8).Value = -4.32223105156559
$ws.Cells.Item(29, 3).Value = 6.261336803436279
$ws.Cells.Item(29, 4).Value = -9.548392295837402
$ws.Cells.Item(29, 5).Value = 26.6091365814209
$ws.Cells.Item(29, 6).Value = 3.495839129919295
$ws.Cells.Item(29, 7).Value = 0.1920398131184147
$ws.Cells.Item(29, 8).Value = -4.956773451004905
$ws.Cells.Item(30, 3).Value = -5.242365837097168
$ws.Cells.Item(30, 4).Value = 18.12157821655273
$ws.Cells.Item(30, 5).Value = 3.32082748413086
$ws.Cells.Item(30, 6).Value = 4.441888874974746
$ws.Cells.Item(30, 7).Value = -4.70699503503997
$ws.Cells.Item(30, 8).Value = -3.40560439948378
$ws.Cells.Item(31, 3).Value = -17.76671600341797
$ws.Cells.Item(31, 4).Value = -38.77998352050781
$ws.Cells.Item(31, 5).Value = 33.62932968139648
$ws.Cells.Item(31, 6).Value = 1.070874737597066
$ws.Cells.Item(31, 7).Value = -0.5835215826144202
$ws.Cells.Item(31, 8).Value = -1.073033490400209
